$d = $word.ActiveDocument

# --- Change 1: merge the leading-spaces run with the "Experienced QA Analyst..." run ---
$d.Content.Find.Execute(
    "     Experienced QA Analyst with a demonstrated history of working in the financial service industry specializing Data governance and Data quality. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "     Experienced QA Analyst with a demonstrated history of working in the financial service industry specializing Data governance and Data quality. ",
    2) | Out-Null

# --- Change 2: merge "A solution oriented..." + "outlook. Proven" + " expertise in" into one run ---
$d.Content.Find.Execute(
    "A solution oriented, dedicated, creative professional with quick grasping ability and problem solving outlook. Proven expertise in",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A solution oriented, dedicated, creative professional with quick grasping ability and problem solving outlook. Proven expertise in",
    2) | Out-Null

# --- Change 3: merge "Quality Assurance " + "role" into one run ---
$d.Content.Find.Execute(
    "Quality Assurance role",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Quality Assurance role",
    2) | Out-Null

# --- Change 4: append a new run " Jenkins  SoapUI" right after the "Solving" run ---
$rng = $d.Content
$rng.Find.Execute("Solving", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" Jenkins  SoapUI")
$rng.Font.Color = 4604471
